$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "4.21.0"
$ws.Range("C9").Value = "102.4.431"
$ws.Range("C12").Value = "1.4.1.986"
$ws.Range("C25").Value = "3.8.5"
$ws.Range("C30").Value = "1.12.0-beta"
